$wb = $excel.ActiveWorkbook

function Set-Row {
    param($ws, $row, $aVal, $bVal, $cVal)
    if ($null -ne $aVal) { $ws.Cells.Item($row, 1).Value = $aVal }
    if ($null -ne $bVal) { $ws.Cells.Item($row, 2).Value = $bVal }
    if ($null -ne $cVal) { $ws.Cells.Item($row, 3).Value = $cVal }
}

$ws = $wb.Worksheets.Item("MANAGER")

Set-Row $ws 1191 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://online.r-m.co.il/login.php' 0.203
Set-Row $ws 1192 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://system.serveandcheck.com/login.php' 2.709
Set-Row $ws 1193 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://eu.checker-soft.com/gfk-ukraine/login.php' 1.12
Set-Row $ws 1194 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://ipsos-russia.com/login.php' 1.244
Set-Row $ws 1195 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://tyaonline.com/login.php' 6.936
Set-Row $ws 1196 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://marketest.checker.co.il/login.php' 1.818
Set-Row $ws 1197 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://ms-online.co.il/login.php' 1.172
Set-Row $ws 1198 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://www.misonline.co.il/login.php' 1.19
Set-Row $ws 1199 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://tovanot.checker.co.il/login.php' 1.11
Set-Row $ws 1200 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://csc.ajis-group.co.jp/login.php' 4.222
Set-Row $ws 1201 '10.13.2022 18:56 (Kyiv+Israel) 15:56 (UTC) 00:56 (Japan) 21:26 (India)' 'https://ru.checker-soft.com/profpoint-ru/login.php' 1.473
Set-Row $ws 1202 '10.13.2022 18:56 (Kyiv+Israel) 15:56 (UTC) 00:56 (Japan) 21:26 (India)' 'https://www.imystery.ru/login.php' 0.961
Set-Row $ws 1203 '10.13.2022 18:56 (Kyiv+Israel) 15:56 (UTC) 00:56 (Japan) 21:26 (India)' 'https://eu.checker-soft.com/testing/login.php' 0.87

Set-Row $ws 1204 $null '*****' $null

$ws = $wb.Worksheets.Item("SHOPPER")

Set-Row $ws 1191 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://online.r-m.co.il/c_login.php' 0.546
Set-Row $ws 1192 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://system.serveandcheck.com/c_login.php' 0.453
Set-Row $ws 1193 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://eu.checker-soft.com/gfk-ukraine/c_login.php' 0.458
Set-Row $ws 1194 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://ipsos-russia.com/c_login.php' 2.608
Set-Row $ws 1195 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://tyaonline.com/c_login.php' 1.678
Set-Row $ws 1196 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://marketest.checker.co.il/c_login.php' 0.351
Set-Row $ws 1197 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://ms-online.co.il/c_login.php' 0.345
Set-Row $ws 1198 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://www.misonline.co.il/c_login.php' 0.447
Set-Row $ws 1199 '10.13.2022 18:55 (Kyiv+Israel) 15:55 (UTC) 00:55 (Japan) 21:25 (India)' 'https://tovanot.checker.co.il/c_login.php' 0.97
Set-Row $ws 1200 '10.13.2022 18:56 (Kyiv+Israel) 15:56 (UTC) 00:56 (Japan) 21:26 (India)' 'https://csc.ajis-group.co.jp/c_login.php' 1.692
Set-Row $ws 1201 '10.13.2022 18:56 (Kyiv+Israel) 15:56 (UTC) 00:56 (Japan) 21:26 (India)' 'https://ru.checker-soft.com/profpoint-ru/c_login.php' 0.362
Set-Row $ws 1202 '10.13.2022 18:56 (Kyiv+Israel) 15:56 (UTC) 00:56 (Japan) 21:26 (India)' 'https://www.imystery.ru/c_login.php' 0.508
Set-Row $ws 1203 '10.13.2022 18:56 (Kyiv+Israel) 15:56 (UTC) 00:56 (Japan) 21:26 (India)' 'https://eu.checker-soft.com/testing/c_login.php' 0.538

Set-Row $ws 1204 $null '*****' $null

